# Rename the worksheet from the location-specific name to the generic
# "Ativos" name (the app now drives the sheet name from the database
# rather than hard-coding a particular unit name), and keep the
# workbook's print area definition in sync with the new sheet name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Ativos"
$ws.PageSetup.PrintArea = "A1:M100"
